# Update countries & provincias Spain
# - Refresh "last updated" timestamp
# - Refresh totals for Estados Unidos (row 4) and China (row 10)
# - Bolivia overtakes Guinea in the ranking (new Bolivia numbers, others shift down a row)
# - Islas Virgenes de los Estados Unidos overtakes Fiyi in the ranking (numbers unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 03:22"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 738830
$ws.Range("C4").Value = 38
$ws.Range("D4").Value = 68285
$ws.Range("E4").Value = 631531

# China (row 10)
$ws.Range("B10").Value = 82735
$ws.Range("C10").Value = 16
$ws.Range("D10").Value = 77062
$ws.Range("E10").Value = 1041

# Bolivia moves ahead of Guinea / Uruguay / Kirguistan with refreshed counts;
# those three keep their own numbers, just shifted one row down.
$ws.Range("A99").Value = "Bolivia"
$ws.Range("B99").Value = 520
$ws.Range("C99").Value = 27
$ws.Range("D99").Value = 31
$ws.Range("E99").Value = 457
$ws.Range("F99").Value = 3
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 32

$ws.Range("A100").Value = "Guinea"
$ws.Range("B100").Value = 518
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 65
$ws.Range("E100").Value = 450
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 3

$ws.Range("A101").Value = "Uruguay"
$ws.Range("B101").Value = 517
$ws.Range("C101").Value = 0
$ws.Range("D101").Value = 298
$ws.Range("E101").Value = 210
$ws.Range("F101").Value = 14
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 9

$ws.Range("A102").Value = "Kirguistan"
$ws.Range("B102").Value = 506
$ws.Range("C102").Value = 0
$ws.Range("D102").Value = 130
$ws.Range("E102").Value = 371
$ws.Range("F102").Value = 5
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 5

# Row 103 (Honduras) is unaffected by the reshuffle.

# Islas Virgenes de los Estados Unidos swaps places with Fiyi (identical totals).
$ws.Range("A183").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("A184").Value = "Fiyi"
